$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns I and J, matching the style of the existing header row (H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I2:J32
$data = @(
    @(7, 8),
    @(7, 8),
    @(9, 9),
    @(9, 9),
    @(6, 8),
    @(7, 7),
    @(7, 8),
    @(8, 8),
    @(7, 7),
    @(7, 7),
    @(8, 8),
    @(6, 6),
    @(7, 7),
    @(8, 8),
    @(5, 6),
    @(9, 9),
    @(6, 6),
    @(7, 7),
    @(7, 7),
    @(10, 10),
    @(6, 7),
    @(6, 6),
    @(8, 8),
    @(6, 6),
    @(6, 6),
    @(7, 7),
    @(6, 6),
    @(5, 5),
    @(5, 5),
    @(3, 3),
    @(3, 3)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
